$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 170.338511186
$ws.Range("C2").Value = 0.0000258227967094

$ws.Range("B3").Value = 170338.511186
$ws.Range("C3").Value = 0.7746839012819999

$ws.Range("B4").Value = 629179.3587677283
$ws.Range("C4").Value = 3.559943665756238

$ws.Range("B5").Value = 12583.58717535457
$ws.Range("C5").Value = 0.07119887331512477
